$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.003662661169252107"
$ws.Range("D3").Value = [double]"2.50620950453622E-19"
$ws.Range("D4").Value = [double]"3.052775686449299E-19"
$ws.Range("D5").Value = [double]"0.002268980687014664"
$ws.Range("D6").Value = [double]"0.002118953267204091"
$ws.Range("D7").Value = [double]"2.481648079184332E-17"
$ws.Range("D8").Value = [double]"2.846335152723769E-17"
$ws.Range("D9").Value = [double]"0.002347931971754905"
$ws.Range("D10").Value = [double]"5.485580910164068E-16"
